$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "project"
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("I1").VerticalAlignment = -4160
$ws.Range("I1").Borders.Item(7).LineStyle = 1
$ws.Range("I1").Borders.Item(10).LineStyle = 1
